# Penalty Reward System (unfinished) - shift forecast weeks forward by one
# and update the MyForecast + Summary figures accordingly.

$wb = $excel.ActiveWorkbook
$wsForecast = $wb.Worksheets.Item("Forecast Comparison")
$wsSummary  = $wb.Worksheets.Item("Summary")

# --- Forecast Comparison sheet -------------------------------------------
# Column B holds week-start dates stored as literal text (not real dates),
# so force a text number format before writing them to avoid Excel
# re-interpreting the strings as date serials.
$wsForecast.Range("B2:B17").NumberFormat = "@"

$weekDates = @{
    2  = "2025-01-12"
    3  = "2025-01-19"
    4  = "2025-01-26"
    5  = "2025-02-02"
    6  = "2025-02-09"
    7  = "2025-02-16"
    8  = "2025-02-23"
    9  = "2025-03-02"
    10 = "2025-03-09"
    11 = "2025-03-16"
    12 = "2025-03-23"
    13 = "2025-03-30"
    14 = "2025-04-06"
    15 = "2025-04-13"
    16 = "2025-04-20"
    17 = "2025-04-27"
}

$myForecast = @{
    2  = 2
    3  = 2
    4  = 2
    5  = 3
    6  = 2
    7  = 2
    8  = 2
    9  = 2
    10 = 2
    11 = 2
    12 = 2
    13 = 2
    14 = 2
    15 = 2
    16 = 3
    17 = 3
}

foreach ($row in 2..17) {
    $wsForecast.Cells.Item($row, 2).Value = $weekDates[$row]
    $wsForecast.Cells.Item($row, 4).Value = $myForecast[$row]
}

# --- Summary sheet ---------------------------------------------------------
# Every value in column B on this sheet is stored as literal text, even the
# numeric-looking ones, so force text formatting before writing.
$wsSummary.Range("B2:B15").NumberFormat = "@"

$wsSummary.Range("B2").Value  = "2022-12-25 to 2025-01-05"
$wsSummary.Range("B5").Value  = "0"
$wsSummary.Range("B9").Value  = "38"
$wsSummary.Range("B10").Value = "18"
$wsSummary.Range("B11").Value = "9"
$wsSummary.Range("B12").Value = "3"
$wsSummary.Range("B13").Value = "2025-04-20"
$wsSummary.Range("B14").Value = "2"
$wsSummary.Range("B15").Value = "2025-01-12"
